$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Medial-Occipital_Visual"
$ws.Range("A3").Value = "Occipital-Lateral(L)_Visual"
$ws.Range("A4").Value = "Medial-Lateral(L)_Visual"
$ws.Range("A5").Value = "Medial-Lateral(R)_Visual"
$ws.Range("A6").Value = "ACC-RPFC(L)_Salience"
$ws.Range("A7").Value = "AInsula(L)-RPFC(R)_Salience"
$ws.Range("A8").Value = "LPFC(L)-PPC(L)_FP"
